# Add the team's name as a new text box on the title slide (slide 1),
# placed between the main title and the "May 6, 2023" subtitle, matching
# the layout/format PowerPoint uses for a single-line, shape-to-text
# autosized text box.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Position/size are specified in points; COM converts to EMU internally
# (EMU / 12700 = points). Target EMU: off (250575, 2642552), ext (3744936, 369332).
$left   = 250575   / 12700.0
$top    = 2642552  / 12700.0
$width  = 3744936  / 12700.0
$height = 369332   / 12700.0

$box = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$box.Name = "TextBox 1"

$tf = $box.TextFrame
$tf.WordWrap = $false
$tf.AutoSize = 1

$box.Fill.Visible = $false

$tr = $tf.TextRange
$tr.Text = "TEAM THEANO CAPSTONE PROJECT"
$tr.Font.Name = "PT Sans"
$tr.Font.Size = 18
$tr.LanguageID = "en-NG"
